$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.480903045476927
$ws.Range("B2").Value = -1.282556239139729
$ws.Range("A3").Value = -0.4855914895904065
$ws.Range("B3").Value = -0.616898906967759
$ws.Range("A4").Value = -0.8648869244303575
$ws.Range("B4").Value = -0.7311905286725847
$ws.Range("A5").Value = -0.7655245086053081
$ws.Range("B5").Value = -0.6893182275130685
$ws.Range("A6").Value = 0.8181679629558203
$ws.Range("B6").Value = 0.642404867836407
$ws.Range("A7").Value = -0.07037208432078354
$ws.Range("B7").Value = 0.04498653625319461
$ws.Range("A8").Value = 0.7822745837455743
$ws.Range("B8").Value = 0.5771386703682702
$ws.Range("A9").Value = 0.3163357981425192
$ws.Range("B9").Value = 0.3009017501146894
$ws.Range("A10").Value = -0.181234671156298
$ws.Range("B10").Value = -0.04247692004725151
